$d = $word.ActiveDocument

# --- 1) Remove the pre-existing _GoBack bookmark (it sat at the very end of
#        the document, right after a manual line break) before we add the
#        new one further up - bookmark names must stay unique.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2) Find the "Ttulo3" heading paragraph that still reads
#        "Analizar y comparar los principales Frameworks de creación de
#        videojuegos 2D para no programadores." and rebuild its runs to
#        match the edited wording / proofing marks / relocated bookmark.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Analizar y comparar los principales Frameworks*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph"
}

$start = $target.Range.Start
$end = $target.Range.End
$r = $d.Range($start, $end)

$r.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:t>Conocer</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t xml:space="preserve"> los principales </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Frameworks</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> de creaci&#243;n de videojuegos 2D para no programadores.</w:t></w:r>
</w:p>
"@)
